# Generate Report for Handback
# The file "502a4e9c-e0f6-48d4-bfa2-6d6278c5b9bc.md" has finished its
# localization round trip: flip its status from "Ready for handoff" to
# "Handed back: in sync with en-US" on every sheet, and stamp the new
# "Latest Handback DateTime" on the locale sheets.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# --- Overview sheet: row for 502a4e9c...md is row 3 (B = zh-cn status, C = de-de status)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status

# --- zh-cn sheet: row for 502a4e9c...md is row 3 (C = Status, H = Latest Handback DateTime)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("H3").Value = "2016-03-20 20:40:01"

# --- de-de sheet: row for 502a4e9c...md is row 3 (C = Status, H = Latest Handback DateTime)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("H3").Value = "2016-03-20 20:40:08"
